$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F2:F12
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 69
$ws1.Range("F3").Value = 543
$ws1.Range("F4").Value = 178
$ws1.Range("F5").Value = 278
$ws1.Range("F6").Value = 392
$ws1.Range("F7").Value = 243
$ws1.Range("F8").Value = 2308
$ws1.Range("F9").Value = 385
$ws1.Range("F10").Value = 5743
$ws1.Range("F11").Value = 139
$ws1.Range("F12").Value = 376

# Sheet "演出" (sheet2) - F2:F4
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 52
$ws2.Range("F3").Value = 11
$ws2.Range("F4").Value = 14

# Sheet "全部类型" (sheet4) - F2:F15
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 69
$ws4.Range("F3").Value = 52
$ws4.Range("F4").Value = 543
$ws4.Range("F5").Value = 178
$ws4.Range("F6").Value = 278
$ws4.Range("F7").Value = 392
$ws4.Range("F8").Value = 243
$ws4.Range("F9").Value = 11
$ws4.Range("F10").Value = 14
$ws4.Range("F11").Value = 2308
$ws4.Range("F12").Value = 385
$ws4.Range("F13").Value = 5743
$ws4.Range("F14").Value = 139
$ws4.Range("F15").Value = 376
